$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.741.50"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.628.33"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'214.27"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'4.25"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.854.19"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "1.625.91"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'0.551"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "0.0₃0760"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "'62.80"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "25.742.31"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'0.999"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'191.34"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'1.81"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "'142.42"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").Value = "'0.123"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'0.0493"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").Value = "'1.58"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").Value = "'2.38"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "1.138.58"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").Value = "'2.50"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "'0.541"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "'0.0156"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "'5.54"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "'100.61"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "1.764.55"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'55.10"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.44"
$ws.Range("E49").Value = "  +5.55%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.417"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  -0.51%  "
